# Applies the "feat: add 2022-Q3 data" edit:
#   1. Inserts a new summary row into the "总计" sheet for "2022-Q3"
#      (35 funds held, 2.8 亿元 market value), pushing the older quarters
#      down by one row.
#   2. Inserts a brand-new worksheet named "2022-Q3" right after "总计"
#      holding the per-fund holdings detail for that quarter.
#
# NOTE: sheet handles returned by Worksheets.Item(<index>) track the
# *position*, not the sheet identity - once the sheet order changes
# (Add/Move/Delete) an old handle silently resolves to whatever now sits
# at that index. So every handle used here is re-fetched *by name*
# immediately before it's used, after any operation that could reorder
# sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" worksheet right after "总计", using
# "2022-Q2" purely as a formatting donor (same header/index-column style).
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

$template = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Item("2022-Q3")
$template.Range("A1:H2").Copy($q3.Range("A1"))
$q3.Range("A2:H2").Copy()
$q3.Range("A3:H36").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Fund code / name / scale / position / ratio / market-value / rank
$q3data = @(
    @('011637','广发沪港深价值成长混合A','12.26','85.19','4.21','0.5161',9),
    @('005598','广发中小盘精选混合A','9.52','93.42','4.69','0.4465',4),
    @('013955','广发中小盘精选混合C','6.86','93.42','4.69','0.3217',4),
    @('009740','博时研究臻选三年持有期灵活配置混合A','7.61','82.29','2.81','0.2138',8),
    @('001521','国寿安保成长优选股票','6.22','91.90','3.14','0.1953',8),
    @('217021','招商优势企业混合','2.91','79.74','5.32','0.1548',5),
    @('398061','中海消费混合','3.91','85.30','3.78','0.1478',10),
    @('501098','建信科技创新 3 年封闭运作灵活配置混合','3.36','77.90','3.26','0.1095',8),
    @('011815','恒越优势精选混合','3.22','90.91','3.16','0.1018',6),
    @('013250','红土创新智能制造混合','1.31','93.66','7.67','0.1005',2),
    @('501076','鹏华创新动力混合（LOF）','5.77','40.77','1.25','0.0721',7),
    @('011845','博时周期优选混合A','2.21','78.30','2.67','0.0590',10),
    @('015071','鑫元专精特新混合A','2.65','74.01','2.00','0.0530',10),
    @('015031','博时远见回报混合C','1.17','76.26','3.23','0.0378',4),
    @('013028','恒越品质生活混合','1.35','89.92','2.74','0.0370',9),
    @('014212','博时研究优享混合A','0.80','79.50','4.12','0.0330',4),
    @('000166','中海信息产业精选混合','0.77','89.31','4.21','0.0324',7),
    @('011340','博时战略新材料主题混合A','0.91','79.55','3.34','0.0304',4),
    @('015030','博时远见回报混合A','0.79','76.26','3.23','0.0255',4),
    @('011341','博时战略新材料主题混合C','0.60','79.55','3.34','0.0200',4),
    @('008082','国寿安保研究精选混合A','0.37','92.74','3.85','0.0142',8),
    @('009741','博时研究臻选三年持有期灵活配置混合C','0.49','82.29','2.81','0.0138',8),
    @('011638','广发沪港深价值成长混合C','0.32','85.19','4.21','0.0135',9),
    @('006072','民生加银创新成长混合A','0.45','93.02','2.97','0.0134',10),
    @('014913','博时研究回报混合A','0.30','78.32','4.18','0.0125',5),
    @('008083','国寿安保研究精选混合C','0.16','92.74','3.85','0.0062',8),
    @('015072','鑫元专精特新混合C','0.25','74.01','2.00','0.0050',10),
    @('014914','博时研究回报混合C','0.10','78.32','4.18','0.0042',5),
    @('014213','博时研究优享混合C','0.07','79.50','4.12','0.0029',4),
    @('011846','博时周期优选混合C','0.10','78.30','2.67','0.0027',10),
    @('012415','德邦上证G60综指增强A','0.09','92.83','2.41','0.0022',6),
    @('004795','富荣福鑫灵活配置混合C','0.05','88.59','2.80','0.0014',10),
    @('004794','富荣福鑫灵活配置混合A','0.01','88.59','2.80','0.0003',10),
    @('012416','德邦上证G60综指增强C','0.01','92.83','2.41','0.0002',6),
    @('014929','民生加银创新成长混合C','0.00','93.02','2.97',0,10)
)

# Columns B-G hold text (fund codes / formatted numbers-as-text), so force
# text format before assigning to avoid Excel's automatic numeric coercion
# (which would e.g. eat the leading zero of a fund code).
$q3.Range("B2:G36").NumberFormat = "@"

$r = 2
foreach ($row in $q3data) {
    $q3.Cells.Item($r, 1).Value = $r - 2
    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 7).Value = $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# Last row's market value (G36) is a genuine 0, stored as a number (not
# text) in the source data - restore numeric type/format for that one cell.
$q3.Range("G36").NumberFormat = "General"
$q3.Cells.Item(36, 7).Value = 0

# ---------------------------------------------------------------------
# Step 2: "总计" (summary) sheet - insert new row 2 for 2022-Q3
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary = $wb.Worksheets.Item("总计")
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 35
$summary.Cells.Item(2, 4).Value = 2.8

# Re-sequence the index column (A) for the remaining rows, which used to
# start at 0 on row 2 and now start on row 3.
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# Restore "2020-Q4" (the last sheet) as the active tab, matching the
# source workbook's view state.
$lastSheet = $wb.Worksheets.Item("2020-Q4")
$lastSheet.Activate()

Write-Output "2022-Q3 sheet added and summary updated"
